$wb = $excel.ActiveWorkbook

# --- Sheet references -------------------------------------------------
$wsCadastro     = $wb.Worksheets.Item("Cadastro")
$wsLancamentos  = $wb.Worksheets.Item("Lançamentos")

# --- 1. Lançamentos: register a new stock entry for "Caneta esferográfica
#        preta" (inserting a real row above the totals row so the table
#        keeps its totals row / structured references intact). ---------
$wsLancamentos.Rows.Item(6).Insert()

$loLancamentos = $wsLancamentos.ListObjects.Item("TblLancamentos")
$loLancamentos.Resize($wsLancamentos.Range("A3:E7"))

$wsLancamentos.Range("A6").Value = "Caneta esferográfica preta"
$wsLancamentos.Range("B6").Value = 43605
$wsLancamentos.Range("C6").Value = 50
$wsLancamentos.Range("D6").Value = 2

# --- 2. Cadastro: calculate SALDO with SOMASE (SUMIF) against the
#        lançamentos table for each product. ----------------------------
$loCadastro = $wsCadastro.ListObjects.Item("TblCadastro")
$colSaldo = $loCadastro.ListColumns.Item("SALDO")
$colSaldo.DataBodyRange.Formula = "=SUMIF(TblLancamentos[PRODUTO],TblCadastro[[#This Row],[PRODUTO]],TblLancamentos[ENTRADA])-SUMIF(TblLancamentos[PRODUTO],TblCadastro[[#This Row],[PRODUTO]],TblLancamentos[SAÍDA])"

# --- 3. Make "Cadastro" the active sheet/tab, with G1 selected on both
#        affected sheets (matches the workbook's new activeTab). -------
$wsLancamentos.Activate()
$wsLancamentos.Range("G1").Select()

$wsCadastro.Activate()
$wsCadastro.Range("G1").Select()
